$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "BP-S-001"
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 6000
$ws.Range("F2").Value = 0

# Update row 3 with new values
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "CL-C-007"
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 8000
$ws.Range("F3").Value = 0

# Delete rows 4 through 7 (old data no longer needed)
$ws.Range("A4:F7").Delete()
